# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newer counts pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 145
    $ws.Range("F3").Value = 222
    $ws.Range("F4").Value = 3762
    $ws.Range("F5").Value = 387
}
